$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has columns A..G = code, name, status, group-name, group-code,
# category-code, category-name (columns D..G = 4..7).
#
# The edit permutes the 4 "codeforiati" columns (D,E,F,G) on every row
# (including the header row) via the cycle:
#   new D = old E
#   new E = old G
#   new F = old D
#   new G = old F
#
# Read everything into memory first (so later writes don't clobber values
# that a later column read still needs), then write the permuted values
# back in a second pass.

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$colD = @{}
$colE = @{}
$colF = @{}
$colG = @{}

for ($r = 1; $r -le $lastRow; $r++) {
    $colD[$r] = $ws.Cells.Item($r, 4).Value2
    $colE[$r] = $ws.Cells.Item($r, 5).Value2
    $colF[$r] = $ws.Cells.Item($r, 6).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $colE[$r]
    $ws.Cells.Item($r, 5).Value = $colG[$r]
    $ws.Cells.Item($r, 6).Value = $colD[$r]
    $ws.Cells.Item($r, 7).Value = $colF[$r]
}
